$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.0253537328125
$ws.Range("D2").Value = 1.033949085430784
$ws.Range("E2").Value = 1.02572052766284
$ws.Range("F2").Value = 1.042551370869524
$ws.Range("I2").Value = 1.031964981092079
$ws.Range("J2").Value = 1.030523225886379
$ws.Range("K2").Value = 1.03674996485774
$ws.Range("L2").Value = 1.028545272226384
$ws.Range("M2").Value = 1.045327738519362
$ws.Range("N2").Value = 1.031986687069813
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.026326344641302
$ws.Range("D3").Value = 1.034704128189619
$ws.Range("E3").Value = 1.026546260452082
$ws.Range("F3").Value = 1.043541893128195
$ws.Range("I3").Value = 1.032145954726745
$ws.Range("J3").Value = 1.031135016161285
$ws.Range("K3").Value = 1.037314238114099
$ws.Range("L3").Value = 1.029178317863239
$ws.Range("M3").Value = 1.046128644787433
$ws.Range("N3").Value = 1.032599346157073
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.026955886464339
$ws.Range("D4").Value = 1.035192455085987
$ws.Range("E4").Value = 1.027081130378522
$ws.Range("F4").Value = 1.044183097746271
$ws.Range("I4").Value = 1.032261183157837
$ws.Range("J4").Value = 1.031530499557018
$ws.Range("K4").Value = 1.03767846369907
$ws.Range("L4").Value = 1.029587861217564
$ws.Range("M4").Value = 1.046646542859849
$ws.Range("N4").Value = 1.032995391184591
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.027220592354674
$ws.Range("D5").Value = 1.03539768993415
$ws.Range("E5").Value = 1.027306123981297
$ws.Range("F5").Value = 1.044452723841226
$ws.Range("I5").Value = 1.032309175998512
$ws.Range("J5").Value = 1.031696667614772
$ws.Range("K5").Value = 1.037831368463583
$ws.Range("L5").Value = 1.029760013354575
$ws.Range("M5").Value = 1.046864184351187
$ws.Range("N5").Value = 1.033161795220047
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.027265040347121
$ws.Range("D6").Value = 1.035432146344481
$ws.Range("E6").Value = 1.027343909239191
$ws.Range("F6").Value = 1.044497998951336
$ws.Range("I6").Value = 1.032317207847964
$ws.Range("J6").Value = 1.031724562476602
$ws.Range("K6").Value = 1.037857029171808
$ws.Range("L6").Value = 1.029788917264631
$ws.Range("M6").Value = 1.046900722426609
$ws.Range("N6").Value = 1.03318972969578
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.026959423295459
$ws.Range("D7").Value = 1.035195197672652
$ws.Range("E7").Value = 1.027084136227818
$ws.Range("F7").Value = 1.044186700254523
$ws.Range("I7").Value = 1.03226182620636
$ws.Range("J7").Value = 1.031532720269809
$ws.Range("K7").Value = 1.037680507668655
$ws.Range("L7").Value = 1.029590161602079
$ws.Range("M7").Value = 1.046649451322981
$ws.Range("N7").Value = 1.032997615051049
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.025682390374104
$ws.Range("D8").Value = 1.034204303854124
$ws.Range("E8").Value = 1.025999470142763
$ws.Range("F8").Value = 1.042886065850981
$ws.Range("I8").Value = 1.03202652953652
$ws.Range("J8").Value = 1.03073006268396
$ws.Range("K8").Value = 1.036940848800624
$ws.Range("L8").Value = 1.028759228911419
$ws.Range("M8").Value = 1.04559847870959
$ws.Range("N8").Value = 1.032193817599361
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.023433631247175
$ws.Range("D9").Value = 1.032456466220865
$ws.Range("E9").Value = 1.024092530804401
$ws.Range("F9").Value = 1.040596300629568
$ws.Range("I9").Value = 1.031597585252105
$ws.Range("J9").Value = 1.029312757289944
$ws.Range("K9").Value = 1.035630646957431
$ws.Range("L9").Value = 1.027294447914076
$ws.Range("M9").Value = 1.043743957225586
$ws.Range("N9").Value = 1.030774499469119
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.021935526376189
$ws.Range("D10").Value = 1.031290130720696
$ws.Range("E10").Value = 1.022824248196453
$ws.Range("F10").Value = 1.039071276551899
$ws.Range("I10").Value = 1.031302030402584
$ws.Range("J10").Value = 1.02836596908024
$ws.Range("K10").Value = 1.034752641936933
$ws.Range("L10").Value = 1.026317592268516
$ws.Range("M10").Value = 1.042505937587393
$ws.Range("N10").Value = 1.02982636671156
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.021287090377803
$ws.Range("D11").Value = 1.030784845525207
$ws.Range("E11").Value = 1.02227579506747
$ws.Range("F11").Value = 1.038411289068446
$ws.Range("I11").Value = 1.031171784270239
$ws.Range("J11").Value = 1.027955553534696
$ws.Range("K11").Value = 1.034371390812312
$ws.Range("L11").Value = 1.025894533501888
$ws.Range("M11").Value = 1.041969476531493
$ws.Range("N11").Value = 1.029415368328865
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.021046270879901
$ws.Range("D12").Value = 1.030597123008096
$ws.Range("E12").Value = 1.022072184592828
$ws.Range("F12").Value = 1.038166195033697
$ws.Range("I12").Value = 1.031123064723906
$ws.Range("J12").Value = 1.02780303999635
$ws.Range("K12").Value = 1.03422961735519
$ws.Range("L12").Value = 1.025737380302213
$ws.Range("M12").Value = 1.041770153228356
$ws.Range("N12").Value = 1.029262638203801
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.021097925730119
$ws.Range("D13").Value = 1.03063739178646
$ws.Range("E13").Value = 1.022115854762149
$ws.Range("F13").Value = 1.038218766046073
$ws.Range("I13").Value = 1.031133530624039
$ws.Range("J13").Value = 1.027835757685981
$ws.Range("K13").Value = 1.034260035461915
$ws.Range("L13").Value = 1.025771090652709
$ws.Range("M13").Value = 1.041812911354302
$ws.Range("N13").Value = 1.029295402356304
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.021267183374928
$ws.Range("D14").Value = 1.030769329073844
$ws.Range("E14").Value = 1.022258962323727
$ws.Range("F14").Value = 1.038391028397879
$ws.Range("I14").Value = 1.03116776403689
$ws.Range("J14").Value = 1.027942948079509
$ws.Range("K14").Value = 1.034359675029591
$ws.Range("L14").Value = 1.025881543375901
$ws.Range("M14").Value = 1.041953001578782
$ws.Range("N14").Value = 1.029402744972486
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.02137147371161
$ws.Range("D15").Value = 1.030850615107979
$ws.Range("E15").Value = 1.02234715014762
$ws.Range("F15").Value = 1.038497172194152
$ws.Range("I15").Value = 1.031188811291204
$ws.Range("J15").Value = 1.028008982769877
$ws.Range("K15").Value = 1.034421045103092
$ws.Range("L15").Value = 1.025949595596809
$ws.Range("M15").Value = 1.042039308159024
$ws.Range("N15").Value = 1.029468873439688
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.021978566312791
$ws.Range("D16").Value = 1.031323659572528
$ws.Range("E16").Value = 1.022860662515291
$ws.Range("F16").Value = 1.039115085368639
$ws.Range("I16").Value = 1.031310626652502
$ws.Range("J16").Value = 1.028393197585136
$ws.Range("K16").Value = 1.034777921862197
$ws.Range("L16").Value = 1.026345667794595
$ws.Range("M16").Value = 1.042541532590298
$ws.Range("N16").Value = 1.029853633884056
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.022359447593961
$ws.Range("D17").Value = 1.031620320580333
$ws.Range("E17").Value = 1.023182969176103
$ws.Range("F17").Value = 1.039502782474617
$ws.Range("I17").Value = 1.031386431140971
$ws.Range("J17").Value = 1.028634085388378
$ws.Range("K17").Value = 1.03500149555781
$ws.Range("L17").Value = 1.026594094188644
$ws.Range("M17").Value = 1.042856460806401
$ws.Range("N17").Value = 1.030094863775608
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.022581633610689
$ws.Range("D18").Value = 1.031793333324649
$ws.Range("E18").Value = 1.023371034754814
$ws.Range("F18").Value = 1.03972895411841
$ws.Range("I18").Value = 1.031430427661047
$ws.Range("J18").Value = 1.028774547704562
$ws.Range("K18").Value = 1.035131799099456
$ws.Range("L18").Value = 1.02673898992773
$ws.Range("M18").Value = 1.043040115382374
$ws.Range("N18").Value = 1.030235525564392
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.022657397427198
$ws.Range("D19").Value = 1.031852321992931
$ws.Range("E19").Value = 1.023435172023881
$ws.Range("F19").Value = 1.039806078624809
$ws.Range("I19").Value = 1.031445392182031
$ws.Range("J19").Value = 1.028822434331625
$ws.Range("K19").Value = 1.035176211708331
$ws.Range("L19").Value = 1.026788394390559
$ws.Range("M19").Value = 1.043102730421493
$ws.Range("N19").Value = 1.030283480195958
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.022318580109
$ws.Range("D20").Value = 1.031588494198517
$ws.Range("E20").Value = 1.023148381519676
$ws.Range("F20").Value = 1.039461182671014
$ws.Range("I20").Value = 1.031378320680361
$ws.Range("J20").Value = 1.028608244897557
$ws.Range("K20").Value = 1.034977518885946
$ws.Range("L20").Value = 1.026567441124132
$ws.Range("M20").Value = 1.042822675884904
$ws.Range("N20").Value = 1.030068986588327
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.021217340150681
$ws.Range("D21").Value = 1.030730477846459
$ws.Range("E21").Value = 1.022216817684497
$ws.Range("F21").Value = 1.038340299911403
$ws.Range("I21").Value = 1.031157692540545
$ws.Range("J21").Value = 1.027911385022468
$ws.Range("K21").Value = 1.034330338065501
$ws.Range("L21").Value = 1.025849018081193
$ws.Range("M21").Value = 1.041911750087409
$ws.Range("N21").Value = 1.029371137092284
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.020525169872273
$ws.Range("D22").Value = 1.030190794813912
$ws.Range("E22").Value = 1.021631740202246
$ws.Range("F22").Value = 1.037635873375293
$ws.Range("I22").Value = 1.031017006075328
$ws.Range("J22").Value = 1.027472854651441
$ws.Range("K22").Value = 1.033922505560656
$ws.Range("L22").Value = 1.025397257564076
$ws.Range("M22").Value = 1.041338681007732
$ws.Range("N22").Value = 1.02893198395783
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.020892081146435
$ws.Range("D23").Value = 1.03047691093001
$ws.Range("E23").Value = 1.021941840439839
$ws.Range("F23").Value = 1.038009272863015
$ws.Range("I23").Value = 1.031091773069857
$ws.Range("J23").Value = 1.027705364347796
$ws.Range("K23").Value = 1.034138792624331
$ws.Range("L23").Value = 1.025636749744593
$ws.Range("M23").Value = 1.041642507128231
$ws.Range("N23").Value = 1.029164823844622
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.022337046290857
$ws.Range("D24").Value = 1.03160287524659
$ws.Range("E24").Value = 1.023164009977955
$ws.Range("F24").Value = 1.039479979725092
$ws.Range("I24").Value = 1.031381986124955
$ws.Range("J24").Value = 1.028619921237377
$ws.Range("K24").Value = 1.034988353230833
$ws.Range("L24").Value = 1.026579484518913
$ws.Range("M24").Value = 1.042837941953251
$ws.Range("N24").Value = 1.030080679509889
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.024014803730177
$ws.Range("D25").Value = 1.032908524981772
$ws.Range("E25").Value = 1.024584994478614
$ws.Range("F25").Value = 1.041188002379321
$ws.Range("I25").Value = 1.031710170973399
$ws.Range("J25").Value = 1.029679506068031
$ws.Range("K25").Value = 1.03597016916267
$ws.Range("L25").Value = 1.027673191544059
$ws.Range("M25").Value = 1.044223694424629
$ws.Range("N25").Value = 1.031141769072538
